$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Localidad 2 label next to the second "Localidad" block (U4) ---
$ws.Range("U4").Value = "Localidad 2"

# --- New style used for data rows 10 and 11: same border/alignment as the
#     existing header style, but font size 10 instead of 11. Build it once
#     on A10, then copy the resulting format (and only the format) onto the
#     whole A10:AW11 block so every cell in both rows shares one style. ---
$ws.Cells.Item(10, 1).Font.Size = 10
$ws.Cells.Item(10, 1).WrapText = $true
$ws.Cells.Item(10, 1).Copy() | Out-Null
$ws.Range("A10:AW11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 10 data: "SECTOR 123" ---
$row10 = @(1, "SECTOR 123", 0, 2, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($c = 1; $c -le 45; $c++) {
    $ws.Cells.Item(10, $c).Value = $row10[$c - 1]
}
$ws.Cells.Item(10, 46).Formula = "=SUM(J10:AS10)"
$ws.Cells.Item(10, 47).Formula = "=SUM(J10:AS10)"
$ws.Cells.Item(10, 48).Formula = "=SUM(J10:AS10)"

# --- Row 11 data: "NUEVO NNN" ---
$row11 = @(2, "NUEVO NNN", 0, 1, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($c = 1; $c -le 45; $c++) {
    $ws.Cells.Item(11, $c).Value = $row11[$c - 1]
}
$ws.Cells.Item(11, 46).Formula = "=SUM(J11:AS11)"
$ws.Cells.Item(11, 47).Formula = "=SUM(J11:AS11)"
$ws.Cells.Item(11, 48).Formula = "=SUM(J11:AS11)"

# --- Selection moved to the newly-filled row 11 ---
$ws.Range("A11:AW11").Select() | Out-Null
